$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16504778743014913"
$ws1.Range("B2").Value = "go_stims-16504778742594943.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778742845232.csv"
$ws1.Range("B4").Value = "go_stims-1650477874286491.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778743004918.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1650477875553556"
$ws2.Range("B2").Value = "OB-16504778745575228.csv"
$ws2.Range("B3").Value = "OB-1650477874978527.csv"
$ws2.Range("B4").Value = "TB-16504778750255194.csv"
$ws2.Range("B5").Value = "ZB-match_3-16504778744034915.csv"
$ws2.Range("B6").Value = "ZB-match_2-1650477874519528.csv"
$ws2.Range("B7").Value = "OB-16504778746884902.csv"
$ws2.Range("B8").Value = "TB-16504778755414913.csv"
$ws2.Range("B9").Value = "TB-1650477875235494.csv"
$ws2.Range("B10").Value = "ZB-match_2-1650477874324489.csv"

# --- Sheet 3: RS_TO (name only changes) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1650477875559491"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1650477875616527"
$ws4.Range("B2").Value = "MM_stims-16504778755844915.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778755604906.csv"
$ws4.Range("B4").Value = "MM_stims-16504778756005278.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778755844915.csv"
$ws4.Range("B6").Value = "MM_stims-1650477875616527.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778756005278.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778756805212"
$ws5.Range("B2").Value = "vSAT_stims-16504778756485274.csv"
$ws5.Range("B3").Value = "SAT_stims-1650477875619491.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778756324918.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778756645167.csv"
